# Weekly fruit/vegetable price update: insert a new week's price record
# for Ajo (Chino / Primera / China) at Terminal Hortofrutícola Agro
# Chillán, pushing the existing data down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 257; all rows 257:381
# shift down to 258:382 (formatting/styles are carried along).
$ws.Rows.Item(257).Insert()

# Populate the newly inserted row 257 with this week's record.
$ws.Range("A257").Value = 7
$ws.Range("B257").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C257").Value = "Ñuble"
$ws.Range("D257").Value = 45016
$ws.Range("E257").Value = 16
$ws.Range("F257").Value = 100112003
$ws.Range("G257").Value = "Ajo"
$ws.Range("H257").Value = "Chino"
$ws.Range("I257").Value = "Primera"
$ws.Range("J257").Value = 60
$ws.Range("K257").Value = 19000
$ws.Range("L257").Value = 20000
$ws.Range("M257").Value = 19500
$ws.Range("N257").Value = "$/malla 10 kilos"
$ws.Range("O257").Value = "China"
$ws.Range("P257").Value = 1950
$ws.Range("Q257").Value = 10
$ws.Range("R257").Value = "Hortaliza"
